# Update "想去人数" (column F) figures across the sheets, per the upstream
# gh-pages data regeneration (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (Exhibition) ----
$ws1 = $wb.Worksheets.Item("展览")
$updates1 = @{
    4  = 8279
    6  = 110
    7  = 7230
    8  = 1143
    9  = 565
    11 = 728
    13 = 165
    17 = 108
    18 = 11893
    20 = 9
    21 = 148
    22 = 2364
    24 = 3343
    25 = 51
    27 = 2811
    28 = 109
    29 = 31
    31 = 3305
    33 = 2419
    34 = 339
    35 = 1671
    36 = 77
    37 = 113
    38 = 5915
    40 = 16
    41 = 166
    43 = 1122
    45 = 1087
    46 = 1555
    47 = 11
    48 = 107
    49 = 1137
}
foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

# ---- Sheet "演出" (Performance) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F27").Value = 1

# ---- Sheet "本地生活" (Local Life) ----
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 283
$ws3.Range("F3").Value = 422

# ---- Sheet "全部类型" (All Types) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$updates4 = @{
    4  = 283
    5  = 422
    8  = 8279
    10 = 110
    11 = 7230
    12 = 7230
    13 = 1143
    14 = 565
    16 = 728
    18 = 165
    20 = 108
    22 = 11894
    25 = 148
    26 = 2364
    27 = 2364
    28 = 3343
    29 = 2811
    30 = 109
    31 = 31
    33 = 3305
    36 = 2419
    37 = 339
    38 = 1671
    39 = 113
    40 = 5915
    44 = 166
    46 = 1122
    48 = 1087
    49 = 1555
    50 = 107
    51 = 1137
}
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
